$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Generate Report for Handback": a handback .xlf for
# c9fa8c4c-258d-409e-bc48-5595cc57495d was processed for both the zh-cn and
# de-de targets (row 7 of each status sheet). Record the new target file,
# the handback xlf that was applied, the handback datetime, and the fact
# that the handback was generated against a stale source revision (Error
# Detail column). Also widen the Error Detail column so the message is
# readable.
# ---------------------------------------------------------------------------

$latestHandbackMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f198a0ccbc47abb856a748aa48dc21b287c57e61/e2e/c9fa8c4c-258d-409e-bc48-5595cc57495d.md"
$handbackMdDisplay = "c9fa8c4c-258d-409e-bc48-5595cc57495d.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/34082ff938418871adefab0dade3626e9894983d/e2e/c9fa8c4c-258d-409e-bc48-5595cc57495d.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f198a0ccbc47abb856a748aa48dc21b287c57e61/e2e/c9fa8c4c-258d-409e-bc48-5595cc57495d.md."

# --- zh-cn sheet ------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Hyperlinks.Add($wsZh.Range("I7"), $latestHandbackMdUrl, "", "", $handbackMdDisplay)
$wsZh.Range("J7").Value = "c9fa8c4c-258d-409e-bc48-5595cc57495d.a01b54325f8d7feedb3cb8226d895c119d18e12c.zh-cn.xlf"
$wsZh.Range("K7").Value = "2016-08-27 14:42:47"
$wsZh.Range("P7").Value = $errorDetail

$wsZh.Columns.Item(16).ColumnWidth = 39.17

# --- de-de sheet --------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Add($wsDe.Range("I7"), $latestHandbackMdUrl, "", "", $handbackMdDisplay)
$wsDe.Range("J7").Value = "c9fa8c4c-258d-409e-bc48-5595cc57495d.a01b54325f8d7feedb3cb8226d895c119d18e12c.de-de.xlf"
$wsDe.Range("K7").Value = "2016-08-27 14:42:53"
$wsDe.Range("P7").Value = $errorDetail

$wsDe.Columns.Item(16).ColumnWidth = 39.17
